$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove trailing cells from rows 11 and 12 ---
$ws.Range("AS11:AV11").Clear() | Out-Null
$ws.Range("AR12:AS12").Clear() | Out-Null

# --- Row 3: add a new sequence AB3:AJ3 = 1..9 ---
$ws.Range("AB3").Value = 1
$ws.Range("AC3").Value = 2
$ws.Range("AD3").Value = 3
$ws.Range("AE3").Value = 4
$ws.Range("AF3").Value = 5
$ws.Range("AG3").Value = 6
$ws.Range("AH3").Value = 7
$ws.Range("AI3").Value = 8
$ws.Range("AJ3").Value = 9

# --- Row 17 (new): D17:F17 = 1,2,3 ---
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 3
$ws.Rows.Item(17).RowHeight = 18

# --- Row 19 (new): D19:E19 = 1,2 and AP19:AQ19 = "a","b" ---
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("AP19").Value = "a"
$ws.Range("AQ19").Value = "b"
$ws.Rows.Item(19).RowHeight = 18

# --- Row 20 (new): AO20:AR20 = 1,2,3,4 ---
$ws.Range("AO20").Value = 1
$ws.Range("AP20").Value = 2
$ws.Range("AQ20").Value = 3
$ws.Range("AR20").Value = 4
$ws.Rows.Item(20).RowHeight = 18

# --- Row 22: add D22 = 1 ---
$ws.Range("D22").Value = 1

# --- New shared strings must be created in this exact order: p2, p1, h, p3 ---
$ws.Range("AS23").Value = "p2"
$ws.Range("AQ23").Value = "p1"
$ws.Range("AP23").Value = "h"
$ws.Range("AR22").Value = "p3"

# --- Row 24 (new): AP24:AS24 = 1,2,3,4 ---
$ws.Range("AP24").Value = 1
$ws.Range("AQ24").Value = 2
$ws.Range("AR24").Value = 3
$ws.Range("AS24").Value = 4
$ws.Rows.Item(24).RowHeight = 18

# --- Row 26: add AO26 = "p1" (reuses existing shared string) ---
$ws.Range("AO26").Value = "p1"

# --- Row 27 (new): AO27 = "p3", AQ27 = "h" (reuse existing shared strings) ---
$ws.Range("AO27").Value = "p3"
$ws.Range("AQ27").Value = "h"
$ws.Rows.Item(27).RowHeight = 18

# --- Row 28 (new): AO28:AQ28 = 3,2,1 ---
$ws.Range("AO28").Value = 3
$ws.Range("AP28").Value = 2
$ws.Range("AQ28").Value = 1
$ws.Rows.Item(28).RowHeight = 18

# --- Update selection to Y20 ---
$ws.Range("Y20").Select() | Out-Null
